# Weekly update: insert two new daily-price records (Primera & Segunda
# quality "Ajo Chino") on 2023-09-08, right after the existing row 15.
# This pushes the former rows 16-52 down to rows 18-54 and extends the
# sheet's used range from A1:R52 to A1:R54.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 16 (rows below shift down, and the
# inserted rows inherit the surrounding row's formatting - in particular
# the date style on column D).
$ws.Rows("16:17").Insert()

# --- New row 16: "Primera" quality ---
$ws.Range("A16").Value = 1
$ws.Range("B16").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C16").Value = "Arica y Parinacota"
$ws.Range("D16").Value = 45177
$ws.Range("E16").Value = 15
$ws.Range("F16").Value = 100112003
$ws.Range("G16").Value = "Ajo"
$ws.Range("H16").Value = "Chino"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 1100
$ws.Range("K16").Value = 24000
$ws.Range("L16").Value = 25000
$ws.Range("M16").Value = 24727
$ws.Range("N16").Value = "$/caja 10 kilos"
$ws.Range("O16").Value = "China"
$ws.Range("P16").Value = 2473
$ws.Range("Q16").Value = 10
$ws.Range("R16").Value = "Hortaliza"

# --- New row 17: "Segunda" quality ---
$ws.Range("A17").Value = 1
$ws.Range("B17").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C17").Value = "Arica y Parinacota"
$ws.Range("D17").Value = 45177
$ws.Range("E17").Value = 15
$ws.Range("F17").Value = 100112003
$ws.Range("G17").Value = "Ajo"
$ws.Range("H17").Value = "Chino"
$ws.Range("I17").Value = "Segunda"
$ws.Range("J17").Value = 1200
$ws.Range("K17").Value = 22000
$ws.Range("L17").Value = 23000
$ws.Range("M17").Value = 22417
$ws.Range("N17").Value = "$/caja 10 kilos"
$ws.Range("O17").Value = "China"
$ws.Range("P17").Value = 2242
$ws.Range("Q17").Value = 10
$ws.Range("R17").Value = "Hortaliza"
